# Realestate Update resale numbers 2025-01-13 13:48
# Appends a new data row (row 22) to the CityResaleNum sheet with the
# resale figures captured on 2025-01-13 13:48:47.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$row = 22

# Columns A and D contain text that Excel would otherwise auto-convert
# (a date-looking string and a zero-padded number). Force them to be
# stored as text, then reset the cell style back to Normal so no new
# cell formatting is left behind on the row.
$ws.Range("A$row").NumberFormat = "@"
$ws.Range("A$row").Value = "2025-01-13"
$ws.Range("A$row").Style = "Normal"

$ws.Range("B$row").Value = "13:48:47"
$ws.Range("C$row").Value = "Monday"

$ws.Range("D$row").NumberFormat = "@"
$ws.Range("D$row").Value = "02"
$ws.Range("D$row").Style = "Normal"

$ws.Range("E$row").Value = 126885
$ws.Range("F$row").Value = 143588
$ws.Range("G$row").Value = 169365
$ws.Range("H$row").Value = 159600
$ws.Range("I$row").Value = -1
$ws.Range("J$row").Value = 142924
$ws.Range("K$row").Value = -1
$ws.Range("L$row").Value = -1
$ws.Range("M$row").Value = 193075
$ws.Range("N$row").Value = 115486
$ws.Range("O$row").Value = 45825
$ws.Range("P$row").Value = 28512
$ws.Range("Q$row").Value = 65384
$ws.Range("R$row").Value = -1
$ws.Range("S$row").Value = 48767
$ws.Range("T$row").Value = -1
